$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the part count for row 8 ("Connector power in") from 1 to 2.
$ws.Range("C8").Value = 2

# Move the active selection to D7 (matches the saved cursor position).
$ws.Range("D7").Select()
